$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 3933
$ws.Range("K3").Value = 3989
$ws.Range("E4").Value = 1002
$ws.Range("K4").Value = 809
$ws.Range("K5").Value = 277
$ws.Range("K6").Value = 4529
$ws.Range("E7").Value = 12528
$ws.Range("K7").Value = 13537

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K3").Value = 41
$ws.Range("K7").Value = 183

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K3").Value = 273
$ws.Range("K5").Value = 23
$ws.Range("K6").Value = 310
$ws.Range("K7").Value = 919

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K3").Value = 96
$ws.Range("K6").Value = 66
$ws.Range("K7").Value = 285

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 150
$ws.Range("K3").Value = 210
$ws.Range("K6").Value = 163
$ws.Range("K7").Value = 559

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 72
$ws.Range("K7").Value = 228

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 123
$ws.Range("K3").Value = 154
$ws.Range("K6").Value = 140
$ws.Range("K7").Value = 460

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 116
$ws.Range("K4").Value = 47
$ws.Range("K5").Value = 31
$ws.Range("K7").Value = 395
$ws.Range("K8").Value = 919
$ws.Range("K10").Value = 78
$ws.Range("K11").Value = 268
$ws.Range("K14").Value = 69
$ws.Range("K18").Value = 93
$ws.Range("K19").Value = 414
$ws.Range("I20").Value = 307
$ws.Range("J20").Value = 268
$ws.Range("K20").Value = 302
$ws.Range("K25").Value = 59
$ws.Range("K26").Value = 22
$ws.Range("K27").Value = 133
$ws.Range("K29").Value = 715
$ws.Range("K31").Value = 148
$ws.Range("K33").Value = 559
$ws.Range("K35").Value = 19
$ws.Range("K36").Value = 171
$ws.Range("K37").Value = 460
$ws.Range("K41").Value = 115
$ws.Range("K42").Value = 480
$ws.Range("K47").Value = 77
$ws.Range("K51").Value = 163
$ws.Range("K52").Value = 369
$ws.Range("K53").Value = 183
$ws.Range("K55").Value = 152
$ws.Range("K60").Value = 88
$ws.Range("E63").Value = 159
$ws.Range("I63").Value = 105
$ws.Range("K63").Value = 53
$ws.Range("K64").Value = 80
$ws.Range("K66").Value = 46
$ws.Range("K67").Value = 526
$ws.Range("K68").Value = 31
$ws.Range("K73").Value = 125
$ws.Range("K76").Value = 193
$ws.Range("J78").Value = 177
$ws.Range("K78").Value = 164
$ws.Range("K79").Value = 350
$ws.Range("K83").Value = 285
$ws.Range("K84").Value = 97
$ws.Range("K85").Value = 610
$ws.Range("K86").Value = 92
$ws.Range("K87").Value = 19
$ws.Range("K94").Value = 168
$ws.Range("K95").Value = 228
$ws.Range("K98").Value = 70
$ws.Range("E101").Value = 12528
$ws.Range("K101").Value = 13537

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 51
$ws.Range("K7").Value = 148

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K6").Value = 155
$ws.Range("K7").Value = 526

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K2").Value = 28
$ws.Range("K7").Value = 97

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K3").Value = 13
$ws.Range("K6").Value = 43

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 205
$ws.Range("K3").Value = 251
$ws.Range("K6").Value = 202
$ws.Range("K7").Value = 715

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 125
$ws.Range("K7").Value = 414

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K6").Value = 107
$ws.Range("K7").Value = 193

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 69

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 115

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 131
$ws.Range("K3").Value = 156
$ws.Range("K6").Value = 172
$ws.Range("K7").Value = 480

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K6").Value = 36
$ws.Range("K7").Value = 78

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J4").Value = 25
$ws.Range("K4").Value = 14
$ws.Range("J7").Value = 177
$ws.Range("K7").Value = 164

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K6").Value = 57
$ws.Range("K7").Value = 152

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 119
$ws.Range("K7").Value = 350

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K3").Value = 25
$ws.Range("K7").Value = 80

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 104
$ws.Range("K3").Value = 91
$ws.Range("I4").Value = 23
$ws.Range("J4").Value = 26
$ws.Range("I7").Value = 307
$ws.Range("J7").Value = 268
$ws.Range("K7").Value = 302

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K6").Value = 22
$ws.Range("K7").Value = 93

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 67
$ws.Range("K3").Value = 49
$ws.Range("K7").Value = 171

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K3").Value = 127
$ws.Range("K6").Value = 96
$ws.Range("K7").Value = 395

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K2").Value = 48
$ws.Range("K6").Value = 71
$ws.Range("K7").Value = 168

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K3").Value = 24
$ws.Range("K7").Value = 59

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 77

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K3").Value = 13
$ws.Range("K7").Value = 70

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("K3").Value = 5
$ws.Range("K6").Value = 22

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("K3").Value = 10
$ws.Range("K6").Value = 46

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 84
$ws.Range("K7").Value = 268

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 19

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 38
$ws.Range("K7").Value = 125

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K2").Value = 33
$ws.Range("K7").Value = 116

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("K3").Value = 9
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 31

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K3").Value = 34
$ws.Range("K7").Value = 133

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K4").Value = 35
$ws.Range("K7").Value = 92

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K2").Value = 45
$ws.Range("K7").Value = 163

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("K3").Value = 8
$ws.Range("K7").Value = 31

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K2").Value = 31
$ws.Range("K7").Value = 88

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 218
$ws.Range("K3").Value = 202
$ws.Range("K6").Value = 141
$ws.Range("K7").Value = 610

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 96
$ws.Range("K7").Value = 369

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("K2").Value = 14
$ws.Range("K3").Value = 11
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 47

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 19
